$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 14 - this shifts the existing row 14 (and
# everything below it) down to row 15, growing the used range from
# A1:R108 to A1:R109, exactly like Excel's own Insert Row command.
$ws.Rows(14).Insert()

# Populate the newly blank row 14 with the new weekly record.
# Columns that simply repeat the constant values used throughout the
# table (market id/name/region, category info, unit, quantity kind,
# classification) are filled in the same way as every other row.
$ws.Range("A14").Value = 6
$ws.Range("B14").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C14").Value = "Metropolitana"
$ws.Range("D14").Value = "2021-09-23"
$ws.Range("E14").Value = 13
$ws.Range("F14").Value = 100112029
$ws.Range("G14").Value = "Orégano"
$ws.Range("H14").Value = "Sin especificar"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 32
$ws.Range("K14").Value = 9000
$ws.Range("L14").Value = 10000
$ws.Range("M14").Value = 9531
$ws.Range("N14").Value = "$/docena de atados"
$ws.Range("O14").Value = "Región Metropolitana"
$ws.Range("P14").Value = 3177
$ws.Range("Q14").Value = 3
$ws.Range("R14").Value = "Hortaliza"
